$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---- 1. Insert two new rows (17 and 18) to expand first table from 6 to 8 rows ----
# This shifts rows 19.. down by 2, matching dimension B4:G38 and merge-cell shifts.
$ws.Range("B17:B18").EntireRow.Insert()

# ---- helper functions: apply a "look" to a range by copying an existing cell
# format (so border/fill definitions are re-used instead of re-created) and then
# tweaking only the font property that actually differs. ----
function Set-VerdanaBodyStyle($rng) {
    # Verdana 10 (non-bold), same border + center/vcenter alignment as the plain body style
    $tmpl = $ws.Range("B11")
    $tmpl.Copy() | Out-Null
    $rng.PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = 0
    $rng.Font.Name = "Verdana"
    $rng.Font.Size = 10
    $rng.Font.Bold = $false
}

function Set-CalibriSmallStyle($rng) {
    # Keep default Calibri (theme minor font) but shrink size to 10pt, preserving the
    # theme/scheme link (only change Size, do NOT touch Font.Name).
    $tmpl = $ws.Range("B23")
    $tmpl.Copy() | Out-Null
    $rng.PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = 0
    $rng.Font.Size = 10
}

function Set-WrapBodyStyle($rng) {
    # Plain body style (border + center/vcenter) plus word-wrap turned on
    $tmpl = $ws.Range("B23")
    $tmpl.Copy() | Out-Null
    $rng.PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = 0
    $rng.WrapText = $true
}

# ---- 2. Table 1 ("Year : 2019-20") body rows 11-18 ----
Set-VerdanaBodyStyle $ws.Range("B11:G11")
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "Menstrual Hygiene management"
$ws.Range("E11").Value = "Dr. Shilpa Tadurwar"

Set-VerdanaBodyStyle $ws.Range("B12:G12")
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = "Shajyog"
$ws.Range("E12").Value = "Shri. Akash Gholap "

Set-VerdanaBodyStyle $ws.Range("B13:G13")
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = "Karate training"
$ws.Range("D13").Value = "13-23 December 2019"
$ws.Range("E13").Value = "Shri. Datta Kadam, Japan Karate Association"

Set-VerdanaBodyStyle $ws.Range("B14:G14")
$ws.Range("B14").Value = 4
$ws.Range("C14").Value = "Yog shibir"

Set-VerdanaBodyStyle $ws.Range("B15:G15")
$ws.Range("B15").Value = 5
$ws.Range("C15").Value = "Karate Shibir phase II"
$ws.Range("D15").Value = "11days - Feb 2020"
$ws.Range("E15").Value = "Shri. Datta Kadam, Japan Karate Association"

Set-VerdanaBodyStyle $ws.Range("B16:G16")
$ws.Range("B16").Value = 6
$ws.Range("C16").Value = "MCM and diseases"
$ws.Range("E16").Value = "Dr. Shilpa Tadurwar"

Set-VerdanaBodyStyle $ws.Range("B17:G17")
$ws.Range("B17").Value = 7
$ws.Range("C17").Value = "महिला सुरक्षाविषयक कायदे"
$ws.Range("E17").Value = "बाभळगाव पोलीस स्टेशन"

Set-VerdanaBodyStyle $ws.Range("B18:G18")
$ws.Range("B18").Value = 8
$ws.Range("C18").Value = "Health Camp    ( blood test )"
$ws.Range("D18").Value = "5th March 2020"
$ws.Range("E18").Value = "Rotary club of India"

# ---- 3. Table 2 ("Year : 2018-19") body rows 23-26 (plain style, values only) ----
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = "Motivation and time management "
$ws.Range("E23").Value = "Shri Raghunath A Kulkarni, PLGP Latur"

$ws.Range("B24").Value = 2
$ws.Range("C24").Value = "Logical Thinking "
$ws.Range("E24").Value = "Smt. V B Swami, V A I T Latur, Shri Avinash Jadhav and Shri Menkudle, MindLabz, Latur"

$ws.Range("B25").Value = 3
$ws.Range("C25").Value = "Pranayam and Yoga"
$ws.Range("E25").Value = "Shri V B Mundhe, Patanjali Yogpeeth Latur"

$ws.Range("B26").Value = 4
$ws.Range("C26").Value = "महिलांवरील अत्याचारासंबंधी कायदे "
$ws.Range("E26").Value = "Shri. P M Makode, P.I, PTC Bhabhalgaon"

# Row 24 E needs wrap-text style (long text), and a taller custom row height
Set-WrapBodyStyle $ws.Range("E24")
$ws.Range("E24").Value = "Smt. V B Swami, V A I T Latur, Shri Avinash Jadhav and Shri Menkudle, MindLabz, Latur"
$ws.Rows.Item(24).RowHeight = 42

# Row 25 gets a slightly taller custom row height too
$ws.Rows.Item(25).RowHeight = 21

# Row 26 C (Marathi text) uses the small Calibri (10pt) look
Set-CalibriSmallStyle $ws.Range("C26")
$ws.Range("C26").Value = "महिलांवरील अत्याचारासंबंधी कायदे "

# ---- 4. Table 2 row 27 (new row, "Verdana body" look - same as table 1) and row 28 (blank) ----
Set-VerdanaBodyStyle $ws.Range("C27")
$ws.Range("B27").Value = 5
$ws.Range("C27").Value = "Menstrual Hygiene management"
Set-VerdanaBodyStyle $ws.Range("E27")
$ws.Range("E27").Value = "Dr. Shilpa Tadurwar, Suman Industries, Latur"
$ws.Range("B28").Value = 6

# ---- 5. Column widths (best-effort; headless engine quantizes fractional widths) ----
$ws.Columns.Item(3).ColumnWidth = 31.7109375
$ws.Columns.Item(4).ColumnWidth = 22.28515625
$ws.Columns.Item(5).ColumnWidth = 44.42578125

# ---- 6. Sheet view: scroll position + final selection ----
$ws.Activate()
$ws.Range("A10").Select()
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("E26").Select()

Write-Host "done"
